$wb = $excel.ActiveWorkbook

# ---- RQ2 sheet: fill in the raw counts that feed the % formulas ----
$ws2 = $wb.Worksheets.Item("RQ2")
$ws2.Range("B3").Value = 83
$ws2.Range("C3").Value = 86
$ws2.Range("E3").Value = 208
$ws2.Range("F3").Value = 10411

$ws2.Range("B4").Value = 63
$ws2.Range("C4").Value = 67
$ws2.Range("E4").Value = 113
$ws2.Range("F4").Value = 2901

$ws2.Activate()
$ws2.Range("C5").Select()

# ---- RQ3 sheet: fill in the raw counts that feed the % formulas ----
$ws3 = $wb.Worksheets.Item("RQ3")
$ws3.Range("B3").Value = 16
$ws3.Range("C3").Value = 16
$ws3.Range("E3").Value = 33
$ws3.Range("F3").Value = 42

$ws3.Range("B4").Value = 12
$ws3.Range("C4").Value = 12
$ws3.Range("E4").Value = 22
$ws3.Range("F4").Value = 28

$ws3.Activate()
$ws3.Range("F3").Select()

# ---- RQ4 sheet: fill in the raw counts that feed the % formulas ----
$ws4 = $wb.Worksheets.Item("RQ4")
$ws4.Range("B3").Value = 16
$ws4.Range("C3").Value = 7
$ws4.Range("E3").Value = 33
$ws4.Range("F3").Value = 11

$ws4.Range("B4").Value = 12
$ws4.Range("C4").Value = 2
$ws4.Range("E4").Value = 22
$ws4.Range("F4").Value = 3

$ws4.Activate()
$ws4.Range("C5").Select()
